# CSU8 - Manter Consulta: sincronizar casos de uso com os protótipos (#45)
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Fluxo principal de "Consultar" - passo 1
# "1. Ator na tela principal Tela T06 seleciona a opção do sistema Consultas."
# ->
# "1. Ator seleciona a opção do sistema "Consultas", através do menu lateral."
# (the red "Tela T06" callout is dropped entirely)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "1. Ator na tela principal Tela T06 seleciona a opção do sistema Consultas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Ator seleciona a opção do sistema “Consultas”, através do menu lateral.",
    2)

# ---------------------------------------------------------------------------
# Edit 2: Fluxo principal de "Adicionar" - passo 1
# "1. Ator seleciona opção de inserção (Tela T09.1)."
# ->
# "1. Ator seleciona opção de inserção (Tela T09.1), clicando em “Adicionar Consulta”."
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("1. Ator seleciona opção de inserção (", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Find.Execute(").", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.MoveStart(1, -1)
$r.InsertBefore(", clicando em “Adicionar Consulta”")

# ---------------------------------------------------------------------------
# Edit 3: Fluxo principal de "Adicionar" - passo 3
# "3. Ator informa os dados e submete para o sistema."
# ->
# "3. Ator informa os dados e submete para o sistema clicando no botão “Salvar”."
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("3. Ator informa os dados e submete para o sistema.", $true, $false, $false, $false, $false, $true, 1, $false, "3. Ator informa os dados e submete para o sistema clicando no botão “Salvar”.", 2)

# ---------------------------------------------------------------------------
# Edit 4: Fluxo principal de "Remover" - passo 1
# "1. Ator seleciona opção de remoção em uma consulta já listada (Tela T09.1)."
# ->
# "...(Tela T09.1), simbolizada pelo ícone de lixeira."
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("1. Ator seleciona opção de remoção em uma consulta já listada (", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Find.Execute(").", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Collapse(0)
$r3.MoveStart(1, -1)
$r3.InsertBefore(", simbolizada pelo ícone de lixeira")

# ---------------------------------------------------------------------------
# Edit 5: Fluxo principal de "Alterar" - passo 1
# "1. Ator seleciona opção de alterar consulta."
# ->
# "1. Ator seleciona opção de alterar consulta, simbolizada pelo ícone de lápis."
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("1. Ator seleciona opção de alterar consulta.", $true, $false, $false, $false, $false, $true, 1, $false, "1. Ator seleciona opção de alterar consulta, simbolizada pelo ícone de lápis.", 2)

# ---------------------------------------------------------------------------
# Edit 6: Fluxo principal de "Alterar" - passo 3
# "3. Ator informa alterações e submete dados para o sistema."
# ->
# "3. Ator informa alterações e submete dados para o sistema, clicando no botão “Salvar”."
# ---------------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute("3. Ator informa alterações e submete dados para o sistema.", $true, $false, $false, $false, $false, $true, 1, $false, "3. Ator informa alterações e submete dados para o sistema, clicando no botão “Salvar”.", 2)
